$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021504201149167
$ws.Cells.Item(2, 4).Value = 1.024148925270816
$ws.Cells.Item(2, 5).Value = 1.030569388043519
$ws.Cells.Item(2, 6).Value = 1.037359412278249
$ws.Cells.Item(2, 9).Value = 1.024956605289207
$ws.Cells.Item(2, 10).Value = 1.026695296140113
$ws.Cells.Item(2, 11).Value = 1.026978275366268
$ws.Cells.Item(2, 12).Value = 1.033380019414963
$ws.Cells.Item(2, 13).Value = 1.040150521338077
$ws.Cells.Item(2, 14).Value = 1.012961808720614
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022822281272698
$ws.Cells.Item(3, 4).Value = 1.025379397710541
$ws.Cells.Item(3, 5).Value = 1.031849051820889
$ws.Cells.Item(3, 6).Value = 1.038989587157702
$ws.Cells.Item(3, 9).Value = 1.025149224544389
$ws.Cells.Item(3, 10).Value = 1.027649105822715
$ws.Cells.Item(3, 11).Value = 1.028014625284432
$ws.Cells.Item(3, 12).Value = 1.034466800191965
$ws.Cells.Item(3, 13).Value = 1.04158831696394
$ws.Cells.Item(3, 14).Value = 1.01327626733418
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.023668898538241
$ws.Cells.Item(4, 4).Value = 1.026169661659227
$ws.Cells.Item(4, 5).Value = 1.032668029784107
$ws.Cells.Item(4, 6).Value = 1.04002626004213
$ws.Cells.Item(4, 9).Value = 1.025263566911149
$ws.Cells.Item(4, 10).Value = 1.028259631090183
$ws.Cells.Item(4, 11).Value = 1.028678743139714
$ws.Cells.Item(4, 12).Value = 1.035160449999638
$ws.Cells.Item(4, 13).Value = 1.042500085201711
$ws.Cells.Item(4, 14).Value = 1.013477530735501
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02402333358339
$ws.Cells.Item(5, 4).Value = 1.026500484597265
$ws.Cells.Item(5, 5).Value = 1.033010182632656
$ws.Cells.Item(5, 6).Value = 1.040457766577329
$ws.Cells.Item(5, 9).Value = 1.025309179131036
$ws.Cells.Item(5, 10).Value = 1.028514718329304
$ws.Cells.Item(5, 11).Value = 1.02895640505907
$ws.Cells.Item(5, 12).Value = 1.035449789146447
$ws.Cells.Item(5, 13).Value = 1.042878978099423
$ws.Cells.Item(5, 14).Value = 1.013561617263582
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.024082758333838
$ws.Cells.Item(6, 4).Value = 1.026555949355851
$ws.Cells.Item(6, 5).Value = 1.033067506378283
$ws.Cells.Item(6, 6).Value = 1.040529966657626
$ws.Cells.Item(6, 9).Value = 1.025316693720234
$ws.Cells.Item(6, 10).Value = 1.028557456514712
$ws.Cells.Item(6, 11).Value = 1.029002936245726
$ws.Cells.Item(6, 12).Value = 1.035498237828821
$ws.Cells.Item(6, 13).Value = 1.042942337960998
$ws.Cells.Item(6, 14).Value = 1.013575705139166
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.023673640319592
$ws.Cells.Item(7, 4).Value = 1.026174087626466
$ws.Cells.Item(7, 5).Value = 1.032672610049736
$ws.Cells.Item(7, 6).Value = 1.040032042743935
$ws.Cells.Item(7, 9).Value = 1.025264186029903
$ws.Cells.Item(7, 10).Value = 1.028263045757094
$ws.Cells.Item(7, 11).Value = 1.028682459275747
$ws.Cells.Item(7, 12).Value = 1.035164325061802
$ws.Cells.Item(7, 13).Value = 1.042505165288538
$ws.Cells.Item(7, 14).Value = 1.013478656358511
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021950961561384
$ws.Cells.Item(8, 4).Value = 1.024566008161525
$ws.Cells.Item(8, 5).Value = 1.031003743426109
$ws.Cells.Item(8, 6).Value = 1.037914120347351
$ws.Cells.Item(8, 9).Value = 1.025023838324381
$ws.Cells.Item(8, 10).Value = 1.027019027697471
$ws.Cells.Item(8, 11).Value = 1.027329864033912
$ws.Cells.Item(8, 12).Value = 1.033749295977271
$ws.Cells.Item(8, 13).Value = 1.040640299824211
$ws.Cells.Item(8, 14).Value = 1.013068542693131
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018866516452797
$ws.Cells.Item(9, 4).Value = 1.021686124266325
$ws.Cells.Item(9, 5).Value = 1.027992656436935
$ws.Cells.Item(9, 6).Value = 1.034041217842882
$ws.Cells.Item(9, 9).Value = 1.024521078118227
$ws.Cells.Item(9, 10).Value = 1.024775217753211
$ws.Cells.Item(9, 11).Value = 1.024896113840952
$ws.Cells.Item(9, 12).Value = 1.031181578379318
$ws.Cells.Item(9, 13).Value = 1.037210198687805
$ws.Cells.Item(9, 14).Value = 1.012328683823678
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016776095443766
$ws.Cells.Item(10, 4).Value = 1.019733924199578
$ws.Cells.Item(10, 5).Value = 1.025936451627
$ws.Cells.Item(10, 6).Value = 1.031361908522255
$ws.Cells.Item(10, 9).Value = 1.024132002400081
$ws.Cells.Item(10, 10).Value = 1.023243489927823
$ws.Cells.Item(10, 11).Value = 1.023238664319254
$ws.Cells.Item(10, 12).Value = 1.029418395655273
$ws.Cells.Item(10, 13).Value = 1.034824160712409
$ws.Cells.Item(10, 14).Value = 1.011823526825483
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015862521700413
$ws.Cells.Item(11, 4).Value = 1.018880658518714
$ws.Cells.Item(11, 5).Value = 1.025034144668085
$ws.Cells.Item(11, 6).Value = 1.030177992004388
$ws.Cells.Item(11, 9).Value = 1.023950579310555
$ws.Cells.Item(11, 10).Value = 1.022571464573794
$ws.Cells.Item(11, 11).Value = 1.022512408572243
$ws.Cells.Item(11, 12).Value = 1.028642382273308
$ws.Cells.Item(11, 13).Value = 1.03376681375554
$ws.Cells.Item(11, 14).Value = 1.011601873891917
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.01552188886969
$ws.Cells.Item(12, 4).Value = 1.018562497159719
$ws.Cells.Item(12, 5).Value = 1.024697157469606
$ws.Cells.Item(12, 6).Value = 1.029734603605361
$ws.Cells.Item(12, 9).Value = 1.023881229157702
$ws.Cells.Item(12, 10).Value = 1.022320501418112
$ws.Cells.Item(12, 11).Value = 1.022241332530618
$ws.Cells.Item(12, 12).Value = 1.028352219551705
$ws.Cells.Item(12, 13).Value = 1.033370379248619
$ws.Cells.Item(12, 14).Value = 1.011519095914727
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015595014642906
$ws.Cells.Item(13, 4).Value = 1.01863079951477
$ws.Cells.Item(13, 5).Value = 1.024769525752941
$ws.Cells.Item(13, 6).Value = 1.029829877142671
$ws.Cells.Item(13, 9).Value = 1.023896194033447
$ws.Cells.Item(13, 10).Value = 1.022374395046917
$ws.Cells.Item(13, 11).Value = 1.022299539061524
$ws.Cells.Item(13, 12).Value = 1.028414547717694
$ws.Cells.Item(13, 13).Value = 1.033455583664188
$ws.Cells.Item(13, 14).Value = 1.011536872397227
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015834391382982
$ws.Cells.Item(14, 4).Value = 1.018854384254972
$ws.Cells.Item(14, 5).Value = 1.025006326739637
$ws.Cells.Item(14, 6).Value = 1.030141415771285
$ws.Cells.Item(14, 9).Value = 1.023944886916482
$ws.Cells.Item(14, 10).Value = 1.022550747423892
$ws.Cells.Item(14, 11).Value = 1.022490028261904
$ws.Cells.Item(14, 12).Value = 1.028618436627577
$ws.Cells.Item(14, 13).Value = 1.033734119944146
$ws.Cells.Item(14, 14).Value = 1.011595040588975
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015981707388314
$ws.Cells.Item(15, 4).Value = 1.018991979662561
$ws.Cells.Item(15, 5).Value = 1.025151984077016
$ws.Cells.Item(15, 6).Value = 1.030332882268125
$ws.Cells.Item(15, 9).Value = 1.023974627789564
$ws.Cells.Item(15, 10).Value = 1.022659225229961
$ws.Cells.Item(15, 11).Value = 1.022607220282578
$ws.Cells.Item(15, 12).Value = 1.02874380432378
$ws.Cells.Item(15, 13).Value = 1.033905244795317
$ws.Cells.Item(15, 14).Value = 1.011630820557557
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016836546850807
$ws.Cells.Item(16, 4).Value = 1.019790382937088
$ws.Cells.Item(16, 5).Value = 1.025996080005843
$ws.Cells.Item(16, 6).Value = 1.031439975588925
$ws.Cells.Item(16, 9).Value = 1.024143768673987
$ws.Cells.Item(16, 10).Value = 1.02328790306459
$ws.Cells.Item(16, 11).Value = 1.023286680883464
$ws.Cells.Item(16, 12).Value = 1.029469630115146
$ws.Cells.Item(16, 13).Value = 1.034893819225002
$ws.Cells.Item(16, 14).Value = 1.011838175076067
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.017370494962699
$ws.Cells.Item(17, 4).Value = 1.020289053748858
$ws.Cells.Item(17, 5).Value = 1.026522334604101
$ws.Cells.Item(17, 6).Value = 1.032128024205126
$ws.Cells.Item(17, 9).Value = 1.024246388234659
$ws.Cells.Item(17, 10).Value = 1.023679889396437
$ws.Cells.Item(17, 11).Value = 1.023710577193188
$ws.Cells.Item(17, 12).Value = 1.029921540962213
$ws.Cells.Item(17, 13).Value = 1.035507414202213
$ws.Cells.Item(17, 14).Value = 1.011967456697795
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017681128593324
$ws.Cells.Item(18, 4).Value = 1.020579154755297
$ws.Cells.Item(18, 5).Value = 1.026828138746843
$ws.Cells.Item(18, 6).Value = 1.032527063249846
$ws.Cells.Item(18, 9).Value = 1.024304996074685
$ws.Cells.Item(18, 10).Value = 1.023907683246321
$ws.Cells.Item(18, 11).Value = 1.023957003693013
$ws.Cells.Item(18, 12).Value = 1.030183924913779
$ws.Cells.Item(18, 13).Value = 1.035862985061343
$ws.Cells.Item(18, 14).Value = 1.012042583634099
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.017786910185954
$ws.Cells.Item(19, 4).Value = 1.020677942672362
$ws.Cells.Item(19, 5).Value = 1.026932215780301
$ws.Cells.Item(19, 6).Value = 1.032662738912586
$ws.Cells.Item(19, 9).Value = 1.024324768539926
$ws.Cells.Item(19, 10).Value = 1.023985212428724
$ws.Cells.Item(19, 11).Value = 1.024040889553023
$ws.Cells.Item(19, 12).Value = 1.030273187112358
$ws.Cells.Item(19, 13).Value = 1.035983832205785
$ws.Cells.Item(19, 14).Value = 1.012068152577839
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.017313291286309
$ws.Cells.Item(20, 4).Value = 1.020235630427093
$ws.Cells.Item(20, 5).Value = 1.026465991785707
$ws.Cells.Item(20, 6).Value = 1.032054440187842
$ws.Cells.Item(20, 9).Value = 1.02423550736676
$ws.Cells.Item(20, 10).Value = 1.023637920562306
$ws.Cells.Item(20, 11).Value = 1.023665182645307
$ws.Cells.Item(20, 12).Value = 1.029873180382051
$ws.Cells.Item(20, 13).Value = 1.03544182253582
$ws.Cells.Item(20, 14).Value = 1.011953615109249
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015763936762828
$ws.Cells.Item(21, 4).Value = 1.018788578029236
$ws.Cells.Item(21, 5).Value = 1.024936645555059
$ws.Cells.Item(21, 6).Value = 1.030049776086969
$ws.Cells.Item(21, 9).Value = 1.02393060235782
$ws.Cells.Item(21, 10).Value = 1.022498853320192
$ws.Cells.Item(21, 11).Value = 1.022433970373112
$ws.Cells.Item(21, 12).Value = 1.02855844960636
$ws.Cells.Item(21, 13).Value = 1.033652200249888
$ws.Cells.Item(21, 14).Value = 1.011577923890342
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014782312811005
$ws.Cells.Item(22, 4).Value = 1.017871683470127
$ws.Cells.Item(22, 5).Value = 1.023964477703705
$ws.Cells.Item(22, 6).Value = 1.028768335378276
$ws.Cells.Item(22, 9).Value = 1.023727538821973
$ws.Cells.Item(22, 10).Value = 1.021774893245232
$ws.Cells.Item(22, 11).Value = 1.021652252382832
$ws.Cells.Item(22, 12).Value = 1.027720718896747
$ws.Cells.Item(22, 13).Value = 1.032505621248595
$ws.Cells.Item(22, 14).Value = 1.011339126004607
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015303409179502
$ws.Cells.Item(23, 4).Value = 1.018358426334714
$ws.Cells.Item(23, 5).Value = 1.024480859497174
$ws.Cells.Item(23, 6).Value = 1.029449665458933
$ws.Cells.Item(23, 9).Value = 1.023836268849303
$ws.Cells.Item(23, 10).Value = 1.02215942441779
$ws.Cells.Item(23, 11).Value = 1.022067385495817
$ws.Cells.Item(23, 12).Value = 1.028165879880017
$ws.Cells.Item(23, 13).Value = 1.033115490052416
$ws.Cells.Item(23, 14).Value = 1.011465965184834
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.017339141664889
$ws.Cells.Item(24, 4).Value = 1.020259772490157
$ws.Cells.Item(24, 5).Value = 1.026491454235982
$ws.Cells.Item(24, 6).Value = 1.032087696705074
$ws.Cells.Item(24, 9).Value = 1.024240427819488
$ws.Cells.Item(24, 10).Value = 1.023656887082364
$ws.Cells.Item(24, 11).Value = 1.023685697035841
$ws.Cells.Item(24, 12).Value = 1.029895036177813
$ws.Cells.Item(24, 13).Value = 1.035471467779435
$ws.Cells.Item(24, 14).Value = 1.011959870394777
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.019669836466353
$ws.Cells.Item(25, 4).Value = 1.02243624100136
$ws.Cells.Item(25, 5).Value = 1.028779580310905
$ws.Cells.Item(25, 6).Value = 1.03505941172934
$ws.Cells.Item(25, 9).Value = 1.02466050375729
$ws.Cells.Item(25, 10).Value = 1.025361531049194
$ws.Cells.Item(25, 11).Value = 1.025531369678476
$ws.Cells.Item(25, 12).Value = 1.031854337471285
$ws.Cells.Item(25, 13).Value = 1.038114270291738
$ws.Cells.Item(25, 14).Value = 1.01252202766477
